# Updated RAD Summary test cases: refresh the "Date" column (B2:B67) on the
# NumberError sheet with a new run's timestamps.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$values = @(
    "Mon Feb 24 22:38:26 EST 2025",
    "Mon Feb 24 22:38:37 EST 2025",
    "Mon Feb 24 22:38:47 EST 2025",
    "Mon Feb 24 22:38:58 EST 2025",
    "Mon Feb 24 22:39:08 EST 2025",
    "Mon Feb 24 22:39:19 EST 2025",
    "Mon Feb 24 22:39:29 EST 2025",
    "Mon Feb 24 22:39:40 EST 2025",
    "Mon Feb 24 22:39:51 EST 2025",
    "Mon Feb 24 22:40:01 EST 2025",
    "Mon Feb 24 22:40:11 EST 2025",
    "Mon Feb 24 22:40:22 EST 2025",
    "Mon Feb 24 22:40:32 EST 2025",
    "Mon Feb 24 22:40:42 EST 2025",
    "Mon Feb 24 22:40:52 EST 2025",
    "Mon Feb 24 22:41:02 EST 2025",
    "Mon Feb 24 22:41:12 EST 2025",
    "Mon Feb 24 22:41:23 EST 2025",
    "Mon Feb 24 22:41:34 EST 2025",
    "Mon Feb 24 22:41:44 EST 2025",
    "Mon Feb 24 22:41:55 EST 2025",
    "Mon Feb 24 22:42:05 EST 2025",
    "Mon Feb 24 22:42:15 EST 2025",
    "Mon Feb 24 22:42:25 EST 2025",
    "Mon Feb 24 22:42:35 EST 2025",
    "Mon Feb 24 22:42:45 EST 2025",
    "Mon Feb 24 22:42:55 EST 2025",
    "Mon Feb 24 22:43:05 EST 2025",
    "Mon Feb 24 22:43:15 EST 2025",
    "Mon Feb 24 22:43:26 EST 2025",
    "Mon Feb 24 22:43:36 EST 2025",
    "Mon Feb 24 22:43:49 EST 2025",
    "Mon Feb 24 22:43:58 EST 2025",
    "Mon Feb 24 22:44:10 EST 2025",
    "Mon Feb 24 22:44:22 EST 2025",
    "Mon Feb 24 22:44:32 EST 2025",
    "Mon Feb 24 22:44:42 EST 2025",
    "Mon Feb 24 22:44:52 EST 2025",
    "Mon Feb 24 22:45:02 EST 2025",
    "Mon Feb 24 22:45:12 EST 2025",
    "Mon Feb 24 22:45:24 EST 2025",
    "Mon Feb 24 22:45:33 EST 2025",
    "Mon Feb 24 22:45:43 EST 2025",
    "Mon Feb 24 22:45:53 EST 2025",
    "Mon Feb 24 22:46:03 EST 2025",
    "Mon Feb 24 22:46:14 EST 2025",
    "Mon Feb 24 22:46:23 EST 2025",
    "Mon Feb 24 22:46:33 EST 2025",
    "Mon Feb 24 22:46:43 EST 2025",
    "Mon Feb 24 22:46:53 EST 2025",
    "Mon Feb 24 22:47:03 EST 2025",
    "Mon Feb 24 22:47:13 EST 2025",
    "Mon Feb 24 22:47:23 EST 2025",
    "Mon Feb 24 22:47:34 EST 2025",
    "Mon Feb 24 22:47:44 EST 2025",
    "Mon Feb 24 22:47:54 EST 2025",
    "Mon Feb 24 22:48:05 EST 2025",
    "Mon Feb 24 22:48:15 EST 2025",
    "Mon Feb 24 22:48:26 EST 2025",
    "Mon Feb 24 22:48:40 EST 2025",
    "Mon Feb 24 22:48:50 EST 2025",
    "Mon Feb 24 22:48:59 EST 2025",
    "Mon Feb 24 22:49:09 EST 2025",
    "Mon Feb 24 22:49:20 EST 2025",
    "Mon Feb 24 22:49:30 EST 2025",
    "Mon Feb 24 22:49:40 EST 2025"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
